$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 8

# Row 4
$ws.Range("B4").Value = "<have>"
$ws.Range("C4").Value = 11

# Row 5
$ws.Range("C5").Value = 16

# Row 7
$ws.Range("C7").Value = 15

# Row 9
$ws.Range("B9").Value = "<bon>"
$ws.Range("C9").Value = 16

# Row 10
$ws.Range("B10").Value = "<of>"

# Row 11
$ws.Range("C11").Value = 12

# Row 13
$ws.Range("C13").Value = 12

# Row 14
$ws.Range("C14").Value = 12

# Row 15
$ws.Range("C15").Value = 18

# Row 17
$ws.Range("C17").Value = 18

# Row 18
$ws.Range("C18").Value = 13
